$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Public exposure sites table has been refreshed with a new set of locations.
# Columns: A=Location, B=Site, C=Exposure period, D=Notes, E=Exist

$rows = @(
    @("Cheltenham",     "Aldi Cheltenham  280/282 Bay Road, Cheltenham VIC 3192",                                                            "29/12/2020 01:30pm-01:45pm", "Case shopped in store",             "new"),
    @("Lakes Entrance",  "Blue Riviera Hire Boats  Marine Parade, Lakes Entrance VIC 3909",                                                    "29/12/2020 11:15am-12:15pm", "Case hired a boat",                 "new"),
    @("Lakes Entrance",  "Central Hotel Lakes Entrance  321 Esplanade, Lakes Entrance VIC 3909",                                               "30/12/2020 5:00pm-6:30pm",   "Case attended outside premises",    "new"),
    @("Lakes Entrance",  "Darcey Annas Beach Cafe Kiosk Gift Shop Gallery  426 Main Beach Walk Surf Life Saving, Lakes Entrance VIC 3909",     "30/12/2020 11:15am-11:20am", "Case picked up takeaway",           "new"),
    @("Lakes Entrance",  "Woolworths Lakes Entrance 371 Esplanade, Lakes Entrance VIC 3909",                                                   "30/12/2020 6:00pm-6:15pm",   "Case shopped in store",             "new"),
    @("Mentone",        "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194",                                                               "29/12/2020 07:30am-08:00am", "Case shopped in store",             "new"),
    @("Mentone",        "Bunnings Mentone  23-27 Nepean Hwy, Mentone VIC 3194",                                                               "31/12/2020 08:00am-08:30am", "Case shopped in store",             "new"),
    @("Mordialloc",     "Woodlands Golf Club  109 White Street Mordialloc VIC 3195",                                                          "28/12/20 2:00pm-5:00pm",     "Case attended course",              "new"),
    @("Mordialloc",     "Woodlands Golf Club  109 White Street Mordialloc VIC 3195",                                                          "28/12/20 2pm-5:00pm",        "Case attended course",              "old"),
    @("Mount Waverley", "IGA  283 Stephensons Road, Mount Waverley VIC 3149",                                                                 "21/12/2020 2:00pm-2:30pm",   "Case shopped for half an hour",     "new"),
    @("Oakleigh",       "Bunnings Oakleigh  1041 Centre Road, Oakleigh South",                                                                "30/12/2020 11:00am-11:30am", "Case shopped for 30 minutes",       "new")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# Columns widened (auto-fit) to match the new (longer) content.
$ws.Columns.Item(1).ColumnWidth = 13.083333333333332
$ws.Columns.Item(2).ColumnWidth = 89.91666666666664
$ws.Columns.Item(3).ColumnWidth = 24.75
$ws.Columns.Item(4).ColumnWidth = 25.083333333333332

$ws.Range("B2").Select() | Out-Null
